$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 425-426, pushing the existing 425-450 rows down to 427-452.
$ws.Range("A425:R426").EntireRow.Insert()

# Row 425 - new "Primera" quality record for the new week (22-Nov-2021 / serial 44516).
$ws.Range("A425").Value = 8
$ws.Range("B425").Value = "Terminal La Palmera de La Serena"
$ws.Range("C425").Value = "Coquimbo"
$ws.Range("D425").Value = 44516
$ws.Range("E425").Value = 4
$ws.Range("F425").Value = 100112023
$ws.Range("G425").Value = "Brócoli"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 2300
$ws.Range("K425").Value = 550
$ws.Range("L425").Value = 600
$ws.Range("M425").Value = 575
$ws.Range("N425").Value = "$/unidad"
$ws.Range("O425").Value = "Provincia del Elquí"
$ws.Range("P425").Value = 575
$ws.Range("Q425").Value = 1
$ws.Range("R425").Value = "Hortaliza"

# Row 426 - new "Segunda" quality record for the same new week.
$ws.Range("A426").Value = 8
$ws.Range("B426").Value = "Terminal La Palmera de La Serena"
$ws.Range("C426").Value = "Coquimbo"
$ws.Range("D426").Value = 44516
$ws.Range("E426").Value = 4
$ws.Range("F426").Value = 100112023
$ws.Range("G426").Value = "Brócoli"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Segunda"
$ws.Range("J426").Value = 1320
$ws.Range("K426").Value = 450
$ws.Range("L426").Value = 500
$ws.Range("M426").Value = 475
$ws.Range("N426").Value = "$/unidad"
$ws.Range("O426").Value = "Provincia del Elquí"
$ws.Range("P426").Value = 475
$ws.Range("Q426").Value = 1
$ws.Range("R426").Value = "Hortaliza"
